{"js": "// Remove the standalone \"Meta description: ...\" paragraph that sits right\n// under the title, and instead append its descriptive sentence (without the\n// \"Meta description: \" label) as a new final, italic paragraph preceded by\n// a new bold \"Play Esqueleto Mariachi Slot for Free | Review\" paragraph,\n// replacing the old italic image-prompt paragraph that used to be last.\n\nconst body = context.document.body;\n\n// --- Step 1: locate & delete the \"Meta description\" paragraph ---\nconst metaResults = body.search(\"Meta description\", { matchCase: false });\nmetaResults.load(\"items\");\nawait context.sync();\n\nif (metaResults.items.length > 0) {\n  const metaParagraph = metaResults.items[0].paragraphs.getFirst();\n  metaParagraph.delete();\n  await context.sync();\n}\n\n// --- Step 2: find the current last paragraph (the italic image-prompt one) ---\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a brand-new bold paragraph right before it. insertHtml is used\n// (rather than insertParagraph) so the new paragraph does NOT inherit the\n// italic direct-formatting of the paragraph it is being inserted next to.\nlastParagraph.insertHtml(\n  \"<p><b>Play Esqueleto Mariachi Slot for Free | Review</b></p>\",\n  \"Before\"\n);\nawait context.sync();\n\n// --- Step 3: re-resolve the (now shifted) final paragraph and swap its text ---\nparagraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst finalParagraph = paragraphs.items[paragraphs.items.length - 1];\nfinalParagraph.insertText(\n  \"Read our review of Esqueleto Mariachi slot game. Play this online casino game for free and win big prizes with its engaging theme and unique special features.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Remove the standalone \"Meta description: ...\" paragraph that sits right\n# under the title, and instead append its descriptive sentence (without the\n# \"Meta description: \" label) as a new final, italic paragraph preceded by\n# a new bold \"Play Esqueleto Mariachi Slot for Free | Review\" paragraph,\n# replacing the old italic image-prompt paragraph that used to be last.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: locate & delete the \"Meta description\" paragraph ---\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Meta description*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- Step 2: insert a new bold paragraph right before the final paragraph ---\n# (the final paragraph is the italic image-prompt paragraph)\n$n = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($n)\n$lastParagraph.Range.InsertParagraphBefore()\n\n# The freshly inserted (still empty) paragraph is now second-to-last.\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n\n# Work off the text-only portion of the range (i.e. exclude the trailing\n# paragraph mark) so the bold/not-italic formatting we apply doesn't bleed\n# into the paragraph mark itself.\n$newTextRange = $newParagraph.Range.Duplicate\n$newTextRange.MoveEnd(1, -1) | Out-Null   # wdCharacter\n$newTextRange.Text = \"Play Esqueleto Mariachi Slot for Free | Review\"\n\n# Re-fetch the range since assigning .Text can change its extent.\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n$newTextRange = $newParagraph.Range.Duplicate\n$newTextRange.MoveEnd(1, -1) | Out-Null\n$newTextRange.Font.Bold = 1\n$newTextRange.Font.Italic = 0\n\n# --- Step 3: swap the text of the (now-shifted) final paragraph ---\n$finalParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$finalTextRange = $finalParagraph.Range.Duplicate\n$finalTextRange.MoveEnd(1, -1) | Out-Null   # keep the paragraph mark/formatting intact\n$finalTextRange.Text = \"Read our review of Esqueleto Mariachi slot game. Play this online casino game for free and win big prizes with its engaging theme and unique special features.\"\n"}
